# Update dates in the Test Summary Report:
#   01-Oct-24 -> 02-Oct-24  (Start date + several TC_00x rows)
#   04-Oct-24 -> 10-Oct-24  (End date)
$d = $word.ActiveDocument

# Replace the "End date" value first (04-Oct-24 -> 10-Oct-24) while the
# original "01-Oct-24" text is still distinct, avoiding any cross-matching.
$d.Content.Find.Execute("04-Oct-24", $false, $false, $false, $false, $false,
                         $true, 1, $false, "10-Oct-24", 2)

# Replace every "Start date" / TC_00x "01-Oct-24" occurrence with 02-Oct-24.
$d.Content.Find.Execute("01-Oct-24", $false, $false, $false, $false, $false,
                         $true, 1, $false, "02-Oct-24", 2)
